$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country ranking swaps (shared-string reorder equivalent) ---
# Brasil overtakes Portugal in total cases: swap the country labels
# in rows 18/19 so row 18 = Brasil (fresh data) and row 19 = Portugal
# (carries the old row-18 numbers, unchanged).
$ws.Range("A18").Value = "Brasil"
$ws.Range("A19").Value = "Portugal"

# Tunez overtakes Republica de Macedonia: same pattern for rows 78/79.
$ws.Range("A78").Value = "Tunez"
$ws.Range("A79").Value = "Republica de Macedonia"

# --- Refreshed "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 20:52"

# --- Updated per-country statistics ---
$ws.Range("B4").Value = 386587
$ws.Range("C4").Value = 19583
$ws.Range("E4").Value = 352996
$ws.Range("F4").Value = 9150
$ws.Range("G4").Value = 1404
$ws.Range("H4").Value = 12275
$ws.Range("B8").Value = 107458
$ws.Range("C8").Value = 4083
$ws.Range("E8").Value = 69394
$ws.Range("G8").Value = 173
$ws.Range("H8").Value = 1983
$ws.Range("B16").Value = 17840
$ws.Range("C16").Value = 1173
$ws.Range("D16").Value = 3935
$ws.Range("E16").Value = 13530
$ws.Range("G16").Value = 52
$ws.Range("H16").Value = 375
$ws.Range("B18").Value = 12610
$ws.Range("C18").Value = 427
$ws.Range("D18").Value = 127
$ws.Range("E18").Value = 11897
$ws.Range("F18").Value = 296
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 586
$ws.Range("B19").Value = 12442
$ws.Range("C19").Value = 712
$ws.Range("D19").Value = 184
$ws.Range("E19").Value = 11913
$ws.Range("F19").Value = 271
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 345
$ws.Range("B21").Value = 9248
$ws.Range("C21").Value = 344
$ws.Range("D21").Value = 770
$ws.Range("E21").Value = 8413
$ws.Range("F21").Value = 149
$ws.Range("G21").Value = 8
$ws.Range("H21").Value = 65
$ws.Range("E71").Value = 348
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 5
$ws.Range("B78").Value = 623
$ws.Range("C78").Value = 27
$ws.Range("D78").Value = 25
$ws.Range("E78").Value = 575
$ws.Range("F78").Value = 39
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 23
$ws.Range("B79").Value = 599
$ws.Range("C79").Value = 29
$ws.Range("D79").Value = 30
$ws.Range("E79").Value = 543
$ws.Range("F79").Value = 15
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = 26
$ws.Range("F90").Value = 15
$ws.Range("D128").Value = 7
$ws.Range("E128").Value = 98
